$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.911.46"
$ws.Range("E2").Value = "  +0.34%  "
$ws.Range("D3").Value = "4.007.69"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'529.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.67%  "
$ws.Range("D6").Value = "'150.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("D7").Value = "'0.690"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +9.92%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("E10").Value = "  -3.72%  "
$ws.Range("E11").Value = "  -4.46%  "
$ws.Range("D12").Value = "'47.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.67%  "
$ws.Range("D13").Value = "'10.62"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.38%  "
$ws.Range("D14").Value = "4.649.07"
$ws.Range("E14").Value = "  -0.53%  "
$ws.Range("D15").Value = "4.005.18"
$ws.Range("E15").Value = "  -1.28%  "
$ws.Range("E16").Value = "  -2.60%  "
$ws.Range("D17").Value = "'20.54"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.23%  "
$ws.Range("E18").Value = "  -1.18%  "
$ws.Range("D19").Value = "'1.19"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.43%  "
$ws.Range("D20").Value = "71.742.18"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").Value = "'426.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.86%  "
$ws.Range("D22").Value = "'97.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.52%  "
$ws.Range("E23").Value = "  -3.67%  "
$ws.Range("E24").Value = "  +3.07%  "
$ws.Range("D25").Value = "'14.38"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.36%  "
$ws.Range("D26").Value = "'11.25"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -7.66%  "
$ws.Range("E27").Value = "  -3.83%  "
$ws.Range("D28").Value = "'5.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("D29").Value = "'36.53"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.94%  "
$ws.Range("D30").Value = "'3.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +22.94%  "
$ws.Range("D31").Value = "'13.33"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.03%  "
$ws.Range("E32").Value = "  -2.27%  "
$ws.Range("D33").Value = "'675.92"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.39%  "
$ws.Range("D34").Value = "'7.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("D35").Value = "'44.40"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.99%  "
$ws.Range("D36").Value = "'65.47"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.69%  "
$ws.Range("D37").Value = "'0.437"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.46%  "
$ws.Range("D38").Value = "'0.152"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("D39").Value = "0.0₃0825"
$ws.Range("E39").Value = "  -9.51%  "
$ws.Range("D40").Value = "'3.39"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.85%  "
$ws.Range("E41").Value = "  -0.09%  "
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("E43").Value = "  -1.89%  "
$ws.Range("D44").Value = "'3.19"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.83%  "
$ws.Range("E45").Value = "  +1.80%  "
$ws.Range("E46").Value = "  -8.27%  "
$ws.Range("D47").Value = "'3.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.78%  "
$ws.Range("E48").Value = "  +2.64%  "
$ws.Range("E49").Value = "  -6.10%  "
$ws.Range("E50").Value = "  -2.87%  "
$ws.Range("D51").Value = "'146.09"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.57%  "
